$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.029.56"
$ws.Range("E2").Value = "  -4.31%  "
$ws.Range("D3").Value = "1.957.26"
$ws.Range("E3").Value = "  -4.21%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.32"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("E6").Value = "  -3.78%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.23"
$ws.Range("E7").Value = "  -8.56%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.372"
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "56.34"
$ws.Range("E10").Value = "  -5.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0802"
$ws.Range("E11").Value = "  +6.28%  "
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("E13").Value = "  -6.48%  "
$ws.Range("E14").Value = "  -7.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.70"
$ws.Range("E15").Value = "  +4.68%  "
$ws.Range("D16").Value = "2.245.58"
$ws.Range("E16").Value = "  -4.09%  "
$ws.Range("E17").Value = "  -3.91%  "
$ws.Range("D18").Value = "1.959.17"
$ws.Range("E18").Value = "  -4.08%  "
$ws.Range("D19").Value = "35.917.88"
$ws.Range("E19").Value = "  -4.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.68"
$ws.Range("E20").Value = "  -3.67%  "
$ws.Range("D21").Value = "0.0₃0851"
$ws.Range("E21").Value = "  -2.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "234.82"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("E23").Value = "  -3.61%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  -5.74%  "
$ws.Range("E26").Value = "  -4.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.73"
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "159.36"
$ws.Range("E28").Value = "  -3.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.71"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("E30").Value = "  +17.30%  "
$ws.Range("E31").Value = "  -2.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.83"
$ws.Range("E32").Value = "  -7.48%  "
$ws.Range("E33").Value = "  -7.00%  "
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("E35").Value = "  -7.94%  "
$ws.Range("E36").Value = "  +3.02%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.25"
$ws.Range("E38").Value = "  -8.79%  "
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.82"
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.05"
$ws.Range("E40").Value = "  +10.47%  "
$ws.Range("E41").Value = "  -4.60%  "
$ws.Range("E42").Value = "  -1.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.84"
$ws.Range("E43").Value = "  -3.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0211"
$ws.Range("E44").Value = "  -3.64%  "
$ws.Range("E45").Value = "  -5.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.78"
$ws.Range("E46").Value = "  -3.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.95"
$ws.Range("E47").Value = "  -5.96%  "
$ws.Range("E48").Value = "  -7.98%  "
$ws.Range("D49").Value = "1.330.26"
$ws.Range("E49").Value = "  -7.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.81"
$ws.Range("E50").Value = "  -4.53%  "
$ws.Range("D51").Value = "2.137.31"
$ws.Range("E51").Value = "  -3.97%  "
